# Applies the "Added final few parts" edit to the Seeed-AQN-BOM sheet:
#  - Removes the "C1, C3" / "302010097" / 2 row (it was merged/folded away)
#  - Adds a new final row "J12" / "320030017" / 1
#  - Leaves the rest of the (sorted) BOM table untouched, letting the rows
#    above the new entry simply shift up by one after the deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find and delete the row whose Designator column holds "C1, C3".
$found = $ws.Range("A1:A39").Find("C1, C3")
if ($found -ne $null) {
    $ws.Rows.Item($found.Row).Delete()
}

# Determine the first empty row after the existing data and append the new part.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "J12"
$ws.Cells.Item($newRow, 2).NumberFormat = "@"
$ws.Cells.Item($newRow, 2).Value = "320030017"
$ws.Cells.Item($newRow, 3).Value = 1

# Reflect the scrolled/selected view recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A30").Select()
